$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (so it lands after "ValidLogin")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "InvalidLogin"

# Populate the new sheet with header + invalid credential sample data
$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("A2").Value = "abcd"
$newSheet.Range("B2").Value = "xyz"

# Match the selection left on the new sheet
$newSheet.Range("B2").Select() | Out-Null
